# Update "想去人数" (F column) counts across all four sheets, matching the
# gh-pages data refresh at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 309
$ws.Range("F5").Value = 1255
$ws.Range("F7").Value = 296
$ws.Range("F10").Value = 6883
$ws.Range("F13").Value = 29
$ws.Range("F14").Value = 7792
$ws.Range("F16").Value = 46
$ws.Range("F17").Value = 4866
$ws.Range("F19").Value = 2278
$ws.Range("F20").Value = 966
$ws.Range("F21").Value = 4536
$ws.Range("F22").Value = 248
$ws.Range("F26").Value = 278
$ws.Range("F29").Value = 2006
$ws.Range("F30").Value = 16
$ws.Range("F31").Value = 224
$ws.Range("F32").Value = 66
$ws.Range("F33").Value = 529
$ws.Range("F35").Value = 1366
$ws.Range("F37").Value = 2092
$ws.Range("F38").Value = 2177

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 79
$ws.Range("F4").Value = 27

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 250
$ws.Range("F3").Value = 1258
$ws.Range("F4").Value = 83

# Sheet 4: 全部类型 (combined view of the other three sheets)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 250
$ws.Range("F4").Value = 1258
$ws.Range("F5").Value = 83
$ws.Range("F7").Value = 309
$ws.Range("F8").Value = 1255
$ws.Range("F9").Value = 79
$ws.Range("F11").Value = 296
$ws.Range("F14").Value = 6883
$ws.Range("F17").Value = 29
$ws.Range("F18").Value = 7792
$ws.Range("F20").Value = 46
$ws.Range("F21").Value = 4866
$ws.Range("F23").Value = 2278
$ws.Range("F24").Value = 966
$ws.Range("F25").Value = 4536
$ws.Range("F26").Value = 248
$ws.Range("F31").Value = 27
$ws.Range("F32").Value = 278
$ws.Range("F35").Value = 2006
$ws.Range("F36").Value = 16
$ws.Range("F37").Value = 224
$ws.Range("F38").Value = 66
$ws.Range("F39").Value = 529
$ws.Range("F42").Value = 1366
$ws.Range("F44").Value = 2092
$ws.Range("F46").Value = 2177
